# Update internal Excel template:
# - Clear the placeholder sample values out of row 8 (columns B:H) on Sheet1,
#   leaving the cell formatting/styles untouched.
# - Move the active cell selection on Sheet1 to H11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the leftover placeholder text from B8:H8 while keeping formatting.
$ws.Range("B8:H8").ClearContents() | Out-Null

# Update the saved selection/active cell for the sheet.
$ws.Range("H11").Select() | Out-Null
